$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing contents (keep formatting/styles intact) so the cells
# get rewritten fresh in the new player order below.
$ws.Cells.ClearContents()

$data = @(
    @("Darius Garland", "PG", "Cleveland Cavaliers"),
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Daniel Gafford", "PF,C", "Dallas Mavericks"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("Mark Williams", "C", "Charlotte Hornets"),
    @("Trey Murphy III", "SF,PF", "New Orleans Pelicans"),
    @("Ziaire Williams", "SG,SF", "Brooklyn Nets"),
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("Jalen Johnson", "SF,PF", "Atlanta Hawks"),
    @("Dereck Lively II", "C", "Dallas Mavericks")
)

# Write header row + column A (player names) first
$ws.Range("A1").Value = "Oyuncu Adı"
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $data[$i][0]
}

# Then column B (positions)
$ws.Range("B1").Value = "Pozisyon"
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $data[$i][1]
}

# Then column C (teams)
$ws.Range("C1").Value = "Takım"
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Range("C$row").Value = $data[$i][2]
}
